$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")

# Row 7 is the "Experimental" property. Its Value cell (B7) was empty and
# now gets the literal text "false". A plain assignment of the bare string
# "false" gets auto-coerced to a Boolean by Excel (same as typing it into a
# cell), so it is entered with a leading quote to force text, then the
# original cell formatting is restored with a formats-only paste from a
# neighboring cell that already carries the correct style.
$meta.Range("B7").Value = "'false"
$meta.Range("B9").Copy()
$meta.Range("B7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 8 is the "Date" property; refresh its generated timestamp value.
$meta.Range("B8").Value = "2025-11-30T13:08:37+00:00"
